$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 119, shifting the existing rows 119:252 down to 120:253
$ws.Rows("119").Insert()

# Populate the new row 119 with its data (same shape/content as the rest of
# the table -- this market's "Zanahoria" observation for date serial 44601,
# i.e. 2022-02-09).
$ws.Range("A119").Value = 7
$ws.Range("B119").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C119").Value = "Ñuble"
$ws.Range("D119").Value = 44601
$ws.Range("E119").Value = 16
$ws.Range("F119").Value = 100114013
$ws.Range("G119").Value = "Zanahoria"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 120
$ws.Range("K119").Value = 6500
$ws.Range("L119").Value = 7000
$ws.Range("M119").Value = 6750
$ws.Range("N119").Value = "$/saco 20 kilos"
$ws.Range("O119").Value = "Provincia de Diguillín"
$ws.Range("P119").Value = 338
$ws.Range("Q119").Value = 20
$ws.Range("R119").Value = "Hortaliza"

Write-Output "done"
